$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (sharedStrings si 6 and si 9) ---
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# --- Type-changing cells (numeric <-> special text placeholder) ---
$ws.Range("C15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("G15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("C16").Value = 4
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C17").Value = 3
$ws.Range("F17").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("D17").Value = 1
$ws.Range("F17").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("E17").Value = 200
$ws.Range("H17").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("C23").Value = 3
$ws.Range("F23").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("D23").Value = 1
$ws.Range("F23").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = 200
$ws.Range("H23").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("C26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("G26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("G26").PasteSpecial(-4122)

$ws.Range("H26").Value = "'***.*"
$ws.Range("A26").Copy()
$ws.Range("H26").PasteSpecial(-4122)

# --- Plain value-only changes ---
$ws.Range("F15").Value = 2
$ws.Range("M15").Value = -7.692307692307
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 10.526315789473
$ws.Range("I16").Value = 183
$ws.Range("J16").Value = 213
$ws.Range("K16").Value = -14.084507042253
$ws.Range("L16").Value = 17.307692307692
$ws.Range("M16").Value = 53.781512605042
$ws.Range("N16").Value = -84.491525423728
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 163
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 8.666666666666
$ws.Range("L17").Value = 17.266187050359
$ws.Range("M17").Value = 94.047619047619
$ws.Range("N17").Value = -37.307692307692
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 350
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -12.5
$ws.Range("I18").Value = 231
$ws.Range("J18").Value = 229
$ws.Range("K18").Value = 0.873362445414
$ws.Range("L18").Value = 15.5
$ws.Range("M18").Value = 9.478672985781
$ws.Range("N18").Value = -90.962441314554
$ws.Range("C19").Value = 36
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 129
$ws.Range("G19").Value = 153
$ws.Range("H19").Value = -15.686274509803
$ws.Range("I19").Value = 1475
$ws.Range("J19").Value = 1537
$ws.Range("K19").Value = -4.033832140533
$ws.Range("L19").Value = 46.329365079365
$ws.Range("M19").Value = 35.197066911090
$ws.Range("N19").Value = -54.545454545454
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -11.111111111111
$ws.Range("I20").Value = 139
$ws.Range("J20").Value = 155
$ws.Range("K20").Value = -10.322580645161
$ws.Range("L20").Value = -8.552631578947
$ws.Range("M20").Value = 78.205128205128
$ws.Range("N20").Value = -95.388188453881
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 27.906976744186
$ws.Range("F21").Value = 188
$ws.Range("G21").Value = 211
$ws.Range("H21").Value = -10.900473933649
$ws.Range("I21").Value = 2206
$ws.Range("J21").Value = 2296
$ws.Range("K21").Value = -3.919860627177
$ws.Range("L21").Value = 32.174955062911
$ws.Range("M21").Value = 38.134001252348
$ws.Range("N21").Value = -78.559626785887
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 41
$ws.Range("K22").Value = -2.439024390243
$ws.Range("L22").Value = 185.714285714286
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 26
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = -18.75
$ws.Range("L23").Value = -13.333333333333
$ws.Range("M23").Value = 8.333333333333
$ws.Range("C24").Value = 57
$ws.Range("D24").Value = 73
$ws.Range("E24").Value = -21.917808219178
$ws.Range("F24").Value = 210
$ws.Range("G24").Value = 349
$ws.Range("H24").Value = -39.828080229226
$ws.Range("I24").Value = 2728
$ws.Range("J24").Value = 3487
$ws.Range("K24").Value = -21.766561514195
$ws.Range("L24").Value = 28.376470588235
$ws.Range("M24").Value = 83.456624075319
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 17.241379310344
$ws.Range("I25").Value = 303
$ws.Range("J25").Value = 328
$ws.Range("K25").Value = -7.621951219512
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.664451827242
$ws.Range("F26").Value = 2
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 77
$ws.Range("J27").Value = 71
$ws.Range("K27").Value = 8.450704225352
$ws.Range("L27").Value = -7.228915662650
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -50

$excel.CutCopyMode = $false

Write-Output "done"